# Auto-generated edit script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.940.85'
$ws.Range("E2").Value = '  +2.59%  '

$ws.Range("D3").Value = '2.612.11'
$ws.Range("E3").Value = '  +0.77%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '568.93'
$ws.Range("E5").Value = '  -0.71%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.90'
$ws.Range("E6").Value = '  -0.87%  '

$ws.Range("E7").Value = '  -0.23%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.601'
$ws.Range("E8").Value = '  -0.31%  '

$ws.Range("D9").Value = '2.633.52'
$ws.Range("E9").Value = '  +1.16%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.71'
$ws.Range("E10").Value = '  +0.55%  '

$ws.Range("E11").Value = '  +1.88%  '

$ws.Range("E12").Value = '  -0.92%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.375'
$ws.Range("E13").Value = '  +8.25%  '

$ws.Range("D14").Value = '3.075.96'
$ws.Range("E14").Value = '  +0.99%  '

$ws.Range("D15").Value = '60.907.85'
$ws.Range("E15").Value = '  +2.47%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '23.56'
$ws.Range("E16").Value = '  +4.21%  '

$ws.Range("E17").Value = '  +1.72%  '

$ws.Range("D18").Value = '2.625.41'
$ws.Range("E18").Value = '  +1.38%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.69'
$ws.Range("E19").Value = '  +3.15%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '349.74'
$ws.Range("E20").Value = '  +3.36%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.00'
$ws.Range("E21").Value = '  +7.07%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.06'
$ws.Range("E22").Value = '  +13.48%  '

$ws.Range("E23").Value = '  +0.22%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.517'
$ws.Range("E24").Value = '  +13.63%  '

$ws.Range("E25").Value = '  -1.39%  '

$ws.Range("E26").Value = '  +0.28%  '

$ws.Range("E27").Value = '  -0.07%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.74'
$ws.Range("E28").Value = '  +6.07%  '

$ws.Range("E29").Value = '  +1.27%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.83'
$ws.Range("E30").Value = '  +8.22%  '

$ws.Range("E31").Value = '  -0.07%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.27'
$ws.Range("E32").Value = '  +2.84%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '161.72'
$ws.Range("E33").Value = '  +1.48%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.53'
$ws.Range("E34").Value = '  +2.24%  '

$ws.Range("E35").Value = '  +4.15%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.966'
$ws.Range("E36").Value = '  +8.96%  '

$ws.Range("E37").Value = '  +4.59%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.59'
$ws.Range("E38").Value = '  +5.91%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '37.70'
$ws.Range("E39").Value = '  +1.52%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.858'
$ws.Range("E40").Value = '  -2.28%  '

$ws.Range("E41").Value = '  +3.37%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '303.75'
$ws.Range("E42").Value = '  +2.66%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '140.81'
$ws.Range("E43").Value = '  +13.02%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.994'
$ws.Range("E44").Value = '  -0.43%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0985'
$ws.Range("E45").Value = '  +0.54%  '

$ws.Range("E46").Value = '  +1.53%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0548'
$ws.Range("E47").Value = '  +1.63%  '

$ws.Range("E48").Value = '  +3.91%  '

$ws.Range("E49").Value = '  +0.60%  '

$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.51'
$ws.Range("E50").Value = '  +5.04%  '

$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.83'
$ws.Range("E51").Value = '  +7.06%  '
